$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cells used to hold two lines of text each ("x \n(independent
# variable)" / "y \n(dependent variable)") wrapped inside a single tall
# row. Split each into its own row instead: row 1 keeps the short label,
# row 2 gets the parenthetical description.
$ws.Range("A1").Value = "x "
$ws.Range("B1").Value = "y"
$ws.Range("A2").Value = "(independent variable)"
$ws.Range("B2").Value = "(dependent variable)"

# Row 2 should look like row 1 (same fill/font/alignment), so copy the
# formatting down.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# Text no longer wraps over a tall single row, so let both rows go back
# to the sheet's default height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

# The header now spans two rows, so freeze the first two rows instead of
# just the first one. Re-select the original active cell afterwards so
# the view selection is unchanged.
$ws.Activate()
$null = $ws.Range("B6").Select()
$excel.ActiveWindow.FreezePanes = $false
$null = $ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("B6").Select()

# Hide the now-unused helper columns C:D with zero width.
$ws.Range("C1:D1").EntireColumn.ColumnWidth = -0.8333333333333334
$ws.Range("C1:D1").EntireColumn.Hidden = $true
